# SAM TODO.xlsx — "updated todo list with additional items"
#
# 1. Row 2 is de-emphasized from a big bold/14pt "banner" row to a normal
#    15pt row, with the bold emphasis moving from the What/Who columns
#    (B2:C2) to the Priority/notes columns (D2:F2).
# 2. Two new TODO rows are appended at the bottom of the table (rows 54/55),
#    matching the formatting of the existing rows above them.
# 3. The current selection moves to B12 (and the sheet is scrolled back to
#    show column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 restyle -------------------------------------------------------
# Columns B/C (What/Who) lose their bold/14pt "header" look...
$ws.Range("B2:C2").Font.Bold = $false
$ws.Range("B2:C2").Font.Size = 11
# ...while columns D/E/F (Priority/Status/notes) pick up bold/11pt emphasis.
$ws.Range("D2:F2").Font.Bold = $true
$ws.Range("D2:F2").Font.Size = 11
# Row was a tall 18.75pt banner row; restore the default row height.
$ws.Rows.Item(2).AutoFit()

# --- New rows 54 & 55 ------------------------------------------------------
# Insert new rows right after the last existing row (53), picking up that
# row's formatting (status/red "Not done" style in A, B style, etc.)
$ws.Rows.Item(54).Insert(-4121, -4163)
$ws.Rows.Item(55).Insert(-4121, -4163)

# Row 54: Ty/Mike - CSP model bug fixes per AOP
$ws.Range("A54").Value = "Not done"
$ws.Range("C54").Value = "Ty/Mike"
$ws.Range("B54").Value = "CSP model bug fixes per AOP"
$ws.Range("E54").Value = "c"

# Row 55: Steve - Loss diagrams renderer updates and update for tech
$ws.Range("A55").Value = "Not done"
$ws.Range("B55").Value = "Loss diagrams renderer updates and update for tech"
$ws.Range("C55").Value = "Steve"
$ws.Range("E55").Value = "c"

# --- Selection / view -------------------------------------------------------
$ws.Range("B12").Select()
